$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 63-64. This pushes the existing rows 63-79 down to
# rows 65-81, preserving all of their data/formatting in place.
$ws.Range("A63:A64").EntireRow.Insert()

# Fill the two newly-inserted rows with this week's new records.

# Row 63: Sandia, "Segunda" quality entry.
$ws.Cells.Item(63, 1).Value = 1
$ws.Cells.Item(63, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(63, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(63, 4).Value = 45215
$ws.Cells.Item(63, 5).Value = 15
$ws.Cells.Item(63, 6).Value = 100112028
$ws.Cells.Item(63, 7).Value = "Sandia"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Segunda"
$ws.Cells.Item(63, 10).Value = 600
$ws.Cells.Item(63, 11).Value = 500
$ws.Cells.Item(63, 12).Value = 550
$ws.Cells.Item(63, 13).Value = 525
$ws.Cells.Item(63, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(63, 15).Value = "Perú"
$ws.Cells.Item(63, 16).Value = 525
$ws.Cells.Item(63, 17).Value = 1
$ws.Cells.Item(63, 18).Value = "Hortaliza"

# Row 64: Sandia, "Tercera" quality entry.
$ws.Cells.Item(64, 1).Value = 1
$ws.Cells.Item(64, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(64, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(64, 4).Value = 45215
$ws.Cells.Item(64, 5).Value = 15
$ws.Cells.Item(64, 6).Value = 100112028
$ws.Cells.Item(64, 7).Value = "Sandia"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Tercera"
$ws.Cells.Item(64, 10).Value = 800
$ws.Cells.Item(64, 11).Value = 480
$ws.Cells.Item(64, 12).Value = 500
$ws.Cells.Item(64, 13).Value = 490
$ws.Cells.Item(64, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(64, 15).Value = "Perú"
$ws.Cells.Item(64, 16).Value = 490
$ws.Cells.Item(64, 17).Value = 1
$ws.Cells.Item(64, 18).Value = "Hortaliza"

# Give the new D63/D64 date cells the same date-number-format style (s="2")
# used throughout column D.
$ws.Cells.Item(63, 4).NumberFormat = $ws.Cells.Item(65, 4).NumberFormat
$ws.Cells.Item(64, 4).NumberFormat = $ws.Cells.Item(65, 4).NumberFormat
